# User stories #36 - #44, #51 - #57 + new requirments
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 45: complexity value changed 3 -> 1 ---
$ws.Range("C45").Value = 1

# --- Drop the old "H" column content (obsolete remark cells on rows 56-58). ---
# Clearing the whole H49:H65 block also lets the row "spans" shrink back to
# 1:7 and the used-range dimension shrink to column G, matching rows that no
# longer reference column H at all.
$ws.Range("H49:H65").ClearContents()

# --- Row 59: the old placeholder note is replaced by a real user story ---
$ws.Range("B59").Value = "pupils get notification for new uploaded assignment"

# --- Row 60 becomes a normal user-story row (previously a stray note row) ---
$ws.Range("A59:E59").Copy()
$ws.Range("A60:E60").PasteSpecial(-4122)
$ws.Range("G59").Copy()
$ws.Range("G60").PasteSpecial(-4122)

$ws.Range("B60").Value = "The list of all of the classes can be viewed"
$ws.Range("C60").Value = 5
$ws.Range("D60").Value = 2
$ws.Range("E60").Value = [DateTime]"2016-11-16"
$ws.Range("G60").Value = "New"

# --- New user-story rows 61-66, formatted like row 59/60 ---
$rows = @(
    @{ Row = 61; Text = "The list of all of the pupils in a specific class can be viewed" },
    @{ Row = 62; Text = "The list of all of the teachers in a specific class can be viewed" },
    @{ Row = 63; Text = "Pupil's details can be viewed by clicking on him" },
    @{ Row = 64; Text = "Teacher's details can be viewed be clicking on him" },
    @{ Row = 65; Text = "The list of all of the teachers can be viewed" },
    @{ Row = 66; Text = "The list of all of the pupils can be viewed" }
)

foreach ($item in $rows) {
    $r = $item.Row

    $ws.Range("A59:E59").Copy()
    $ws.Range("A$r" + ":E$r").PasteSpecial(-4122)
    $ws.Range("G59").Copy()
    $ws.Range("G$r").PasteSpecial(-4122)

    $ws.Cells.Item($r, 1).Value = $r - 1
    $ws.Range("B$r").Value = $item.Text
    $ws.Range("C$r").Value = 5
    $ws.Range("D$r").Value = 2
    $ws.Range("E$r").Value = [DateTime]"2016-11-16"
    $ws.Range("G$r").Value = "New"
}

# --- Trailing empty ID-only rows 67-71 ---
for ($r = 67; $r -le 71; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
    $ws.Cells.Item($r, 1).Style = $ws.Cells.Item($r - 1, 1).Style
}

# --- View state: scroll position + active selection ---
$excel.ActiveWindow.ScrollRow = 49
$ws.Range("C69").Select()
